$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Refresh the "datetimeFigureOut" date placeholder text (5/22/2025 -> 5/30/2025)
#    on the slide master and every slide layout, mirroring PowerPoint's
#    whole-deck recalculation of date fields on save.
# ---------------------------------------------------------------------------
function Update-DatePlaceholder($shapes, $newText) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.Name -like "Date Placeholder*") {
            $sh.TextFrame.TextRange.Text = $newText
        }
    }
}

Update-DatePlaceholder $p.SlideMaster.Shapes "5/30/2025"
for ($li = 1; $li -le $p.SlideMaster.CustomLayouts.Count; $li++) {
    $layout = $p.SlideMaster.CustomLayouts.Item($li)
    Update-DatePlaceholder $layout.Shapes "5/30/2025"
}

# ---------------------------------------------------------------------------
# 2) Insert two new "Comparaison Résultats" slides before the final
#    "Conclusion" slide (previously slide 7, now pushed to slide 9).
# ---------------------------------------------------------------------------

# --- Slide 7: Comparaison Résultats – Dataset 1 ----------------------------
$s7 = $p.Slides.Add(7, 6)
$s7.Shapes.Item(1).TextFrame.TextRange.Text = "Comparaison Résultats – Dataset 1"

$tb7 = $s7.Shapes.AddTextbox(1, 72, 108, 576, 360)
$tb7.Name = "TextBox 2"
$tb7.Fill.Visible = $false
$tb7.TextFrame.WordWrap = $false

$tr7 = $tb7.TextFrame.TextRange
$tr7.Text = "🔹 Choix de Sienna :" + [char]10 + "- Coût total : 498.76 €" + [char]10 + "- Profit total : 196.61 €" + [char]10 + "- Rendement : 39.42%" + [char]10 + "🔹 Algorithme optimisé :" + [char]10 + "- Coût total : 499.96 €" + [char]10 + "- Profit total : 198.55 €" + [char]10 + "- Rendement : 39.71%" + [char]10 + "📌 Actions communes sélectionnées : 0"
$tr7.Font.Size = 18
$tr7.Font.Name = "Calibri"
$tr7.Font.Color.RGB = 0

$tb7.TextFrame.AutoSize = 1
$tb7.Width = 576
$tb7.Height = 360

# --- Slide 8: Comparaison Résultats – Dataset 2 ----------------------------
$s8 = $p.Slides.Add(8, 6)
$s8.Shapes.Item(1).TextFrame.TextRange.Text = "Comparaison Résultats – Dataset 2"

$tb8 = $s8.Shapes.AddTextbox(1, 72, 108, 576, 360)
$tb8.Name = "TextBox 2"
$tb8.Fill.Visible = $false
$tb8.TextFrame.WordWrap = $false

$tr8 = $tb8.TextFrame.TextRange
$tr8.Text = "🔹 Choix de Sienna :" + [char]10 + "- Coût total : 489.24 €" + [char]10 + "- Profit total : 193.78 €" + [char]10 + "- Rendement : 39.61%" + [char]10 + "🔹 Algorithme optimisé :" + [char]10 + "- Coût total : 499.92 €" + [char]10 + "- Profit total : 197.96 €" + [char]10 + "- Rendement : 39.60%" + [char]10 + "📌 Actions communes sélectionnées : 18"
$tr8.Font.Size = 18
$tr8.Font.Name = "Calibri"
$tr8.Font.Color.RGB = 0

$tb8.TextFrame.AutoSize = 1
$tb8.Width = 576
$tb8.Height = 360

Write-Output "done"
